# Generate Report for Handback
#
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: zh-cn / de-de status cells move from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: fill in "Latest Target File" / "Latest
#    Handback File" / "Latest Handback DateTime" for both data rows, and
#    hyperlink the newly-populated "Latest Target File" cells to the same
#    source-markdown pages linked from column A.
#  - Widen the columns that now hold the longer handback status text / file
#    names so the report stays readable.

$wb = $excel.ActiveWorkbook

$githubBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/0f3538a01032c57efbe4ceed90e1fd71f2a1f2e0/e2e/"
$mdA = "758ea1ca-7f3b-4bef-bdb4-b374c798eb73.md"
$mdB = "881106fb-bcc8-4666-9e6c-6b0829a2a3e0.md"

# ---------------------------------------------------------------------------
# Overview sheet: status text for both locales, on both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de status columns so the longer text fits.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# Per-locale detail sheets (zh-cn, de-de): same shape of edit on each.
# ---------------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; XlfSuffix = "zh-cn.xlf"; HandbackTime = "2016-08-13 15:16:42" },
    @{ Name = "de-de"; XlfSuffix = "de-de.xlf"; HandbackTime = "2016-08-13 15:16:52" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Row 2 -> 758ea1ca...
    $ws.Range("K2").Value = $locale.HandbackTime
    $ws.Range("J2").Value = "758ea1ca-7f3b-4bef-bdb4-b374c798eb73.4d617b7815db51686dd8b758611f752839ba92d6." + $locale.XlfSuffix
    $ws.Hyperlinks.Add($ws.Range("I2"), $githubBase + $mdA, "", "", $mdA) | Out-Null

    # Row 3 -> 881106fb...
    $ws.Range("K3").Value = $locale.HandbackTime
    $ws.Range("J3").Value = "881106fb-bcc8-4666-9e6c-6b0829a2a3e0.5b72af00043d5e072f8bd1fbac7fd69cfc22ea69." + $locale.XlfSuffix
    $ws.Hyperlinks.Add($ws.Range("I3"), $githubBase + $mdB, "", "", $mdB) | Out-Null

    # Widen Status / Latest Target File / Latest Handback File columns.
    $ws.Columns.Item(3).ColumnWidth = 29.1
    $ws.Columns.Item(9).ColumnWidth = 39.1
    $ws.Columns.Item(10).ColumnWidth = 39.1
}
